# Populate the weekly project plan task list and format it the way the
# sheet ends up looking after the edit: three task rows starting at B6,
# a wrapped/empty helper cell next to the second task, and the two
# columns sized to fit their content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Make draw.io diagram of product"
$ws.Range("B7").Value = "Pick suitable micro controller"
$ws.Range("B8").Value = "Make schematic outline frame "

# C7 stays empty but carries a wrap-text style (cellXf with alignment
# wrapText="1") just like in the target sheet.
$ws.Range("C7").WrapText = $true

# Column widths as they ended up in the saved workbook.
$ws.Columns.Item(2).ColumnWidth = 34.33203125
$ws.Columns.Item(3).ColumnWidth = 14.109375

# Final selection left on the sheet.
$null = $ws.Range("C11").Select()
